# Ticket 526: add pricing_pricing_interest_rate_type field
# Adds a new column Q ("pricing_interest_rate_type") with per-row values
# to the "invalid" sheet, plus a trailing empty column R, and updates the
# selection/scroll position on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalid")

# Header (Q1) + per-row values (Q2:Q11) - written with the workbook's
# default (unstyled) cell format, so no explicit style index is stamped.
$ws.Cells.Item(1, 17).Value = "pricing_interest_rate_type"
$ws.Cells.Item(1, 17).Style = "Normal"

$values = @{
    2  = 1
    3  = 2
    4  = 3
    5  = 4
    6  = 5
    7  = 6
    8  = 999
    9  = 0
    10 = 10
    11 = 1000
}
foreach ($row in $values.Keys) {
    $cell = $ws.Cells.Item($row, 17)
    $cell.Value = $values[$row]
    $cell.Style = "Normal"
}

# Trailing empty column R (R1:R11) - create the cells with no explicit
# style/content so they serialize as bare <c r="Rn"/> placeholders.
for ($row = 1; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 18)
    $cell.Value = 0
    $cell.ClearContents()
    $cell.Style = "Normal"
}

# Update the view: scroll window and active selection.
$ws.Activate()
$ws.Range("P15").Select()
